# Append incident-log rows 231-235 to the bottom of the sheet.
# Source data (author's commit appended 5 new rows of shop-floor incident
# records dated 2024-06-12 after the existing data which ended at row 230).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (Bloque, Incidencia, Fecha, Hora, Turno,
# Hora de Reparacion, Tiempo de Reparacion, MTBF).
$newRows = @(
    @("WV50 FILTER", "Core enganchado",            "2024-06-12", "09:15:03", "Mañana", "09:15:05", "0:00:02", "-0.01 minutos"),
    @("WC47 NACP",   "No pone tornillo",            "2024-06-12", "09:31:39", "Mañana", "09:31:40", "0:00:01", "-0.00 minutos"),
    @("WC47 NACP",   "Fallo cámara visión",         "2024-06-12", "09:31:44", "Mañana", "09:31:45", "0:00:01", "0.02 minutos"),
    @("WC47 NACP",   "Fallo en paletizador",        "2024-06-12", "10:37:10", "Mañana", "10:37:11", "0:00:01", "-0.00 minutos"),
    @("WC47 NACP",   "Palet atascado en la curva",  "2024-06-12", "10:37:13", "Mañana", "10:37:13", "0:00:00", "0.01 minutos")
)

$startRow = 231

# Force column C (Fecha) to be stored as plain text rather than letting
# Excel auto-coerce the "yyyy-mm-dd"-shaped string into a date serial
# number, matching the source workbook where every cell (including dates
# and times) is plain text.
$ws.Range("C$startRow`:C235").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
